$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 449.76923
$ws.Range("I55").Value = 466.66666
$ws.Range("J55").Value = 435.2857
$ws.Range("K55").Value = 466.66666
$ws.Range("L55").Value = 435.2857
$ws.Range("M55").Value = -252.66666
$ws.Range("N55").Value = -863.2857
$ws.Range("H69").Value = 7943.478
$ws.Range("J69").Value = 8077.273
$ws.Range("L69").Value = 24231.819
$ws.Range("N69").Value = -25979.819
$ws.Range("H72").Value = 7943.478
$ws.Range("J72").Value = 8077.273
$ws.Range("L72").Value = 72695.45699999999
$ws.Range("N72").Value = -81431.45699999999
$ws.Range("H115").Value = 1207.5714
$ws.Range("I115").Value = 1069.6154
$ws.Range("K115").Value = 3208.8462
$ws.Range("M115").Value = -1641.8462
$ws.Range("H138").Value = 3362.776
$ws.Range("J138").Value = 3581.851
$ws.Range("L138").Value = 10745.553
$ws.Range("N138").Value = -21025.553
$ws.Range("H139").Value = 100327
$ws.Range("I139").Value = 80354.5
$ws.Range("K139").Value = 80354.5
$ws.Range("M139").Value = -75214.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1715031.5
$ws.Range("I2").Value = 2694131.5
$ws.Range("J2").Value = 1606.3334
$ws.Range("K2").Value = 2694131.5
$ws.Range("L2").Value = 1606.3334
$ws.Range("M2").Value = -2694018.5
$ws.Range("N2").Value = -1832.3334
$ws.Range("H32").Value = 8996.275
$ws.Range("I32").Value = 4930.407
$ws.Range("J32").Value = 23633.4
$ws.Range("K32").Value = 4930.407
$ws.Range("L32").Value = 23633.4
$ws.Range("M32").Value = -4643.407
$ws.Range("N32").Value = -24207.4
$ws.Range("H45").Value = 10279901
$ws.Range("I45").Value = 17983928
$ws.Range("K45").Value = 17983928
$ws.Range("M45").Value = -17983551
$ws.Range("H49").Value = 19333.334
$ws.Range("J49").Value = 19333.334
$ws.Range("L49").Value = 19333.334
$ws.Range("N49").Value = -19853.334
$ws.Range("H74").Value = 46164.617
$ws.Range("I74").Value = 2622.6667
$ws.Range("K74").Value = 2622.6667
$ws.Range("M74").Value = -1748.6667
$ws.Range("H77").Value = 46164.617
$ws.Range("I77").Value = 2622.6667
$ws.Range("K77").Value = 13113.3335
$ws.Range("M77").Value = -8745.333500000001
$ws.Range("H110").Value = 993293.1
$ws.Range("I110").Value = 1030048.44
$ws.Range("K110").Value = 1030048.44
$ws.Range("M110").Value = -1028003.44
$ws.Range("H116").Value = 1715031.5
$ws.Range("I116").Value = 2694131.5
$ws.Range("J116").Value = 1606.3334
$ws.Range("K116").Value = 2694131.5
$ws.Range("L116").Value = 1606.3334
$ws.Range("M116").Value = -2691837.5
$ws.Range("N116").Value = -6194.3334
$ws.Range("H132").Value = 4762.273
$ws.Range("I132").Value = 3347.5
$ws.Range("J132").Value = 5570.7144
$ws.Range("K132").Value = 10042.5
$ws.Range("L132").Value = 16712.1432
$ws.Range("M132").Value = -7512.5
$ws.Range("N132").Value = -21772.1432

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1715031.5
$ws.Range("I3").Value = 2694131.5
$ws.Range("J3").Value = 1606.3334
$ws.Range("K3").Value = 2694131.5
$ws.Range("L3").Value = 1606.3334
$ws.Range("M3").Value = -2694017.5
$ws.Range("N3").Value = -1834.3334
$ws.Range("H15").Value = 9989
$ws.Range("J15").Value = 9989
$ws.Range("L15").Value = 9989
$ws.Range("N15").Value = -10443
$ws.Range("H22").Value = 2837.5557
$ws.Range("I22").Value = 3105.4285
$ws.Range("K22").Value = 3105.4285
$ws.Range("M22").Value = -2932.4285
$ws.Range("H80").Value = 493.1
$ws.Range("I80").Value = 548.1818
$ws.Range("J80").Value = 461.21054
$ws.Range("K80").Value = 548.1818
$ws.Range("L80").Value = 461.21054
$ws.Range("M80").Value = 449.8182
$ws.Range("N80").Value = -2457.21054
$ws.Range("H83").Value = 493.1
$ws.Range("I83").Value = 548.1818
$ws.Range("J83").Value = 461.21054
$ws.Range("K83").Value = 2740.909
$ws.Range("L83").Value = 2306.0527
$ws.Range("M83").Value = 2251.091
$ws.Range("N83").Value = -12290.0527
$ws.Range("H99").Value = 5105405
$ws.Range("I99").Value = 6805745
$ws.Range("K99").Value = 6805745
$ws.Range("M99").Value = -6804247
$ws.Range("H134").Value = 3519.4
$ws.Range("I134").Value = 1652.7693
$ws.Range("K134").Value = 4958.3079
$ws.Range("M134").Value = -2423.3079

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6588.846
$ws.Range("I58").Value = 7931.25
$ws.Range("J58").Value = 4441
$ws.Range("K58").Value = 7931.25
$ws.Range("L58").Value = 4441
$ws.Range("M58").Value = -7728.25
$ws.Range("N58").Value = -4847
$ws.Range("H94").Value = 758.1429000000001
$ws.Range("I94").Value = 237.5
$ws.Range("J94").Value = 966.4
$ws.Range("K94").Value = 237.5
$ws.Range("L94").Value = 966.4
$ws.Range("M94").Value = 213.5
$ws.Range("N94").Value = -1868.4
$ws.Range("H105").Value = 1518.8462
$ws.Range("I105").Value = 923
$ws.Range("K105").Value = 923
$ws.Range("M105").Value = 824
$ws.Range("H132").Value = 175659.73
$ws.Range("I132").Value = 128171.125
$ws.Range("K132").Value = 384513.375
$ws.Range("M132").Value = -381983.375
$ws.Range("H136").Value = 6588.846
$ws.Range("I136").Value = 7931.25
$ws.Range("J136").Value = 4441
$ws.Range("K136").Value = 23793.75
$ws.Range("L136").Value = 13323
$ws.Range("M136").Value = -21243.75
$ws.Range("N136").Value = -18423

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 56200
$ws.Range("J37").Value = 56200
$ws.Range("L37").Value = 168600
$ws.Range("N37").Value = -168824
$ws.Range("H50").Value = 5500
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H53").Value = 5500
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H98").Value = 1684.9286
$ws.Range("J98").Value = 1854.7273
$ws.Range("L98").Value = 5564.1819
$ws.Range("N98").Value = -8560.1819
$ws.Range("H106").Value = 5982.609
$ws.Range("J106").Value = 5982.609
$ws.Range("L106").Value = 17947.827
$ws.Range("N106").Value = -19839.827
$ws.Range("H129").Value = 1451.15
$ws.Range("J129").Value = 1991
$ws.Range("L129").Value = 5973
$ws.Range("N129").Value = -15973
$ws.Range("H131").Value = 20840504
$ws.Range("J131").Value = 20841066
$ws.Range("L131").Value = 62523198
$ws.Range("N131").Value = -62533278
$ws.Range("H132").Value = 1512.1111
$ws.Range("I132").Value = 1201
$ws.Range("K132").Value = 10809
$ws.Range("M132").Value = -8279

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 19427202
$ws.Range("I80").Value = 27606132
$ws.Range("J80").Value = 2241
$ws.Range("K80").Value = 27606132
$ws.Range("L80").Value = 2241
$ws.Range("M80").Value = -27605134
$ws.Range("N80").Value = -4237
$ws.Range("H83").Value = 19427202
$ws.Range("I83").Value = 27606132
$ws.Range("J83").Value = 2241
$ws.Range("K83").Value = 138030660
$ws.Range("L83").Value = 11205
$ws.Range("M83").Value = -138025668
$ws.Range("N83").Value = -21189
$ws.Range("H113").Value = 4764556
$ws.Range("I113").Value = 7938068.5
$ws.Range("J113").Value = 4287.143
$ws.Range("K113").Value = 7938068.5
$ws.Range("L113").Value = 4287.143
$ws.Range("M113").Value = -7935898.5
$ws.Range("N113").Value = -8627.143
$ws.Range("H132").Value = 4615
$ws.Range("I132").Value = 4375.875
$ws.Range("J132").Value = 4997.6
$ws.Range("K132").Value = 13127.625
$ws.Range("L132").Value = 14992.8
$ws.Range("M132").Value = -10597.625
$ws.Range("N132").Value = -20052.8
$ws.Range("H139").Value = 94402.125
$ws.Range("J139").Value = 94402.125
$ws.Range("L139").Value = 94402.125
$ws.Range("N139").Value = -104682.125

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 82487.45
$ws.Range("I22").Value = 178472.4
$ws.Range("K22").Value = 178472.4
$ws.Range("M22").Value = -178177.4
$ws.Range("H27").Value = 82487.45
$ws.Range("I27").Value = 178472.4
$ws.Range("K27").Value = 178472.4
$ws.Range("M27").Value = -178365.4
$ws.Range("H46").Value = 6305.278
$ws.Range("J46").Value = 6382.0586
$ws.Range("L46").Value = 6382.0586
$ws.Range("N46").Value = -6758.0586
$ws.Range("H48").Value = 30000
$ws.Range("I48").Value = 20000
$ws.Range("K48").Value = 20000
$ws.Range("M48").Value = -19339

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H132").Value = 100102010
$ws.Range("I132").Value = 142874880
$ws.Range("J132").Value = 298627.34
$ws.Range("K132").Value = 428624640
$ws.Range("L132").Value = 895882.02
$ws.Range("M132").Value = -428622110
$ws.Range("N132").Value = -900942.02
$ws.Range("H136").Value = 5402.2915
$ws.Range("I136").Value = 5197.875
$ws.Range("K136").Value = 15593.625
$ws.Range("M136").Value = -13043.625
$ws.Range("H138").Value = 85489.25
$ws.Range("J138").Value = 85489.25
$ws.Range("L138").Value = 85489.25
$ws.Range("N138").Value = -95769.25
